# Update the crypto price ("D" column) and 1h volume change ("E" column)
# figures for rows 2-51 on the active sheet, per the scheduled GitHub
# Actions refresh of the cryptos list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new Price (or $null if unchanged), new Volume(1h) text.
$updates = @(
    @{ Row = 2;  D = '28.540.88';  E = '  -0.02%  ' },
    @{ Row = 3;  D = '1.565.39';   E = '  -1.67%  ' },
    @{ Row = 4;  D = $null;        E = '  +0.17%  ' },
    @{ Row = 5;  D = '211.75';     E = '  -1.42%  ' },
    @{ Row = 6;  D = $null;        E = '  -0.79%  ' },
    @{ Row = 7;  D = $null;        E = '  +0.21%  ' },
    @{ Row = 8;  D = '46.35';      E = '  +5.59%  ' },
    @{ Row = 9;  D = '24.13';      E = '  +0.42%  ' },
    @{ Row = 10; D = $null;        E = '  -1.87%  ' },
    @{ Row = 11; D = $null;        E = '  -1.69%  ' },
    @{ Row = 12; D = $null;        E = '  -0.66%  ' },
    @{ Row = 13; D = '1.789.32';   E = '  -1.62%  ' },
    @{ Row = 14; D = '1.549.58';   E = '  -2.63%  ' },
    @{ Row = 15; D = '0.519';      E = '  -2.25%  ' },
    @{ Row = 16; D = '28.538.44';  E = '  +0.07%  ' },
    @{ Row = 17; D = $null;        E = '  -3.08%  ' },
    @{ Row = 18; D = '61.95';      E = '  -3.27%  ' },
    @{ Row = 19; D = '227.08';     E = '  -2.72%  ' },
    @{ Row = 20; D = $null;        E = '  -2.38%  ' },
    @{ Row = 21; D = '7.31';       E = '  -2.81%  ' },
    @{ Row = 22; D = $null;        E = '  +0.08%  ' },
    @{ Row = 23; D = $null;        E = '  -6.83%  ' },
    @{ Row = 24; D = '9.13';       E = '  -3.22%  ' },
    @{ Row = 25; D = '2.07';       E = '  +5.94%  ' },
    @{ Row = 26; D = '150.85';     E = '  -0.64%  ' },
    @{ Row = 27; D = $null;        E = '  -2.54%  ' },
    @{ Row = 28; D = $null;        E = '  -3.06%  ' },
    @{ Row = 29; D = $null;        E = '  -3.53%  ' },
    @{ Row = 30; D = $null;        E = '  +0.17%  ' },
    @{ Row = 31; D = $null;        E = '  -1.94%  ' },
    @{ Row = 32; D = $null;        E = '  -3.79%  ' },
    @{ Row = 33; D = $null;        E = '  -1.80%  ' },
    @{ Row = 34; D = $null;        E = '  -0.15%  ' },
    @{ Row = 35; D = '1.392.02';   E = '  -1.92%  ' },
    @{ Row = 36; D = '1.04';       E = '  -0.98%  ' },
    @{ Row = 37; D = $null;        E = '  -3.98%  ' },
    @{ Row = 38; D = '2.35';       E = '  +1.00%  ' },
    @{ Row = 39; D = '2.59';       E = '  +1.24%  ' },
    @{ Row = 40; D = $null;        E = '  -1.21%  ' },
    @{ Row = 41; D = $null;        E = '  -1.70%  ' },
    @{ Row = 42; D = $null;        E = '  +0.17%  ' },
    @{ Row = 43; D = $null;        E = '  -3.54%  ' },
    @{ Row = 44; D = $null;        E = '  +1.31%  ' },
    @{ Row = 45; D = $null;        E = '  -4.28%  ' },
    @{ Row = 46; D = $null;        E = '  -0.19%  ' },
    @{ Row = 47; D = '62.57';      E = '  -3.32%  ' },
    @{ Row = 48; D = '1.701.53';   E = '  -1.69%  ' },
    @{ Row = 49; D = '86.14';      E = '  -2.12%  ' },
    @{ Row = 50; D = '0.0₆0102';   E = '  -1.78%  ' },
    @{ Row = 51; D = $null;        E = '  -0.78%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Range("D$($u.Row)")
        # Several prices (e.g. "211.75") parse as plain numbers, which would
        # make Excel store them as numeric cells instead of text and lose the
        # original two-decimal-group formatting. Force text storage with a
        # leading apostrophe, then restore the default "Normal" style so no
        # stray number-format style lingers on the cell.
        if ($u.D -match '^-?[0-9]+(\.[0-9]+)?$') {
            $cell.Value = "'" + $u.D
            $cell.Style = "Normal"
        } else {
            $cell.Value = $u.D
        }
    }
    $ws.Range("E$($u.Row)").Value = $u.E
}
